$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'Meridian Round Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(2, 3).Value = '₹5,939'
$ws.Cells.Item(3, 2).Value = 'Adele Rectangular Engineered Wood Coffee Table In Classic Walnut Finish'
$ws.Cells.Item(3, 3).Value = '₹2,927'
$ws.Cells.Item(4, 2).Value = 'Awdry Rectangular Engineered Wood Coffee Table In Sonoma Oak Finish'
$ws.Cells.Item(4, 3).Value = '₹2,903'
$ws.Cells.Item(5, 2).Value = 'Altura Rectangular Solid Wood Coffee Table In Two Tone Finish'
$ws.Cells.Item(5, 3).Value = '₹14,453'
$ws.Cells.Item(6, 2).Value = 'Claire Rectangular Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(6, 3).Value = '₹12,725'
$ws.Cells.Item(7, 2).Value = 'Tate Square Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(7, 3).Value = '₹16,554'
$ws.Cells.Item(8, 2).Value = 'Striado Rectangular Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(8, 3).Value = '₹10,947'
$ws.Cells.Item(9, 2).Value = 'Marcel Rectangular Metal Coffee Table In White Gloss Finish'
$ws.Cells.Item(9, 3).Value = '₹11,967'
$ws.Cells.Item(10, 2).Value = 'Renesme Rectangular Solid Wood Coffee Table In Mahogany Finish'
$ws.Cells.Item(10, 3).Value = '₹15,317'
$ws.Cells.Item(11, 2).Value = 'Dyson Abstract Metal Coffee Table In Teak Finish'
$ws.Cells.Item(11, 3).Value = '₹7,679'
$ws.Cells.Item(12, 2).Value = 'Ivara Rectangular Solid Wood Coffee Table In Natural Finish'
$ws.Cells.Item(12, 3).Value = '₹16,049'
$ws.Cells.Item(13, 2).Value = 'Botwin Rectangular Solid Wood Coffee Table In Mahogany Finish'
$ws.Cells.Item(13, 3).Value = '₹9,647'
$ws.Cells.Item(14, 2).Value = 'Zephyr Rectangular Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(14, 3).Value = '₹14,104'
$ws.Cells.Item(15, 2).Value = 'Fring Engineered Wood Side Table In Matte Finish'
$ws.Cells.Item(15, 3).Value = '₹2,399'
$ws.Cells.Item(16, 2).Value = 'Claire Rectangular Solid Wood Coffee Table In Mahogany Finish'
$ws.Cells.Item(16, 3).Value = '₹12,725'
$ws.Cells.Item(17, 2).Value = 'Botwin Rectangular Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(17, 3).Value = '₹9,647'
$ws.Cells.Item(18, 2).Value = 'Epsilon Rectangular Solid Wood Coffee Table In Mahogany Finish'
$ws.Cells.Item(18, 3).Value = '₹11,384'
$ws.Cells.Item(19, 2).Value = 'Dyson Rectangular Metal Coffee Table In Walnut Finish'
$ws.Cells.Item(19, 3).Value = '₹10,529'
$ws.Cells.Item(20, 2).Value = 'Gustowe Rectangular Engineered Wood Coffee Table In Matte Finish'
$ws.Cells.Item(20, 3).Value = '₹2,279'
$ws.Cells.Item(21, 2).Value = 'Striado Rectangular Solid Wood Coffee Table In Mahogany Finish'
$ws.Cells.Item(21, 3).Value = '₹10,947'
$ws.Cells.Item(22, 2).Value = 'Osiris Rectangular Stone Coffee Table In Finish'
$ws.Cells.Item(22, 3).Value = '₹15,677'
$ws.Cells.Item(23, 2).Value = 'Altura Rectangular Solid Wood Coffee Table In Two Tone Finish'
$ws.Cells.Item(23, 3).Value = '₹8,374'
$ws.Cells.Item(24, 2).Value = 'Sylvie Rectangular Solid Wood Coffee Table In Natural Finish'
$ws.Cells.Item(24, 3).Value = '₹11,839'
$ws.Cells.Item(25, 2).Value = 'Florence Oval Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(25, 3).Value = '₹10,223'
$ws.Cells.Item(26, 2).Value = 'Liam Rectangular Engineered Wood Coffee Table In Dark Wenge Finish'
$ws.Cells.Item(26, 3).Value = '₹3,817'
$ws.Cells.Item(27, 2).Value = 'Reeves Rectangular Engineered Wood Coffee Table In Rustic Walnut Finish'
$ws.Cells.Item(27, 3).Value = '₹5,543'
$ws.Cells.Item(28, 2).Value = 'Nitara Oval Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(28, 3).Value = '₹12,095'
$ws.Cells.Item(29, 2).Value = 'Renesme Rectangular Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(29, 3).Value = '₹15,317'
$ws.Cells.Item(30, 2).Value = 'Odette Square Solid Wood Coffee Table In Honey Oak Finish'
$ws.Cells.Item(30, 3).Value = '₹5,919'
$ws.Cells.Item(31, 2).Value = 'Epsilon Rectangular Solid Wood Coffee Table In Teak Finish'
$ws.Cells.Item(31, 3).Value = '₹11,384'
